# Apply updated crypto price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.162.11'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.800.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4723'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +25.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3727'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.22'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07705'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.147'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.71'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.002'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.391'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.413'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.797.87'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001099'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06756'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.433'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.154.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.405'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.393'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.004.57'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.90'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.260'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.045'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09688'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.946'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02385'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2219'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06363'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.67%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6727'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.274'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.239'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.500'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.112'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6172'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.856'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.066'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.184'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07124'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.82%  '
